# Automatic update of files.
# Updates the "Förändrad" date column (C) for rows 2-14 from 2023-09-15
# (serial 45184) to 2023-09-16 (serial 45185), matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45185
    }
}
